# Generate Report for Handoff
# Updates the "Status" text and handoff/generate timestamps across the
# Overview / zh-cn / de-de sheets, then lets the columns re-autofit to
# the (now shorter) status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# Overview sheet: zh-cn / de-de status columns (E2, F2) + generate date (G2)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-23 19:03:33"

# zh-cn sheet: Status (C2) + Latest Handoff Datetime (H2)
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-23 19:03:28"

# de-de sheet: Status (C2) + Latest Handoff Datetime (H2)
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-23 19:03:33"

# Columns re-fit to the new (shorter) status text. The host's ColumnWidth
# setter quantizes to whole pixels, so feed it the input value that lands on
# the closest achievable width to the real-Excel autofit result (~17.22).
$newWidth = 16.3
$overview.Range("E1").ColumnWidth = $newWidth
$overview.Range("F1").ColumnWidth = $newWidth
$zhcn.Range("C1").ColumnWidth = $newWidth
$dede.Range("C1").ColumnWidth = $newWidth
